$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================================
# This script reproduces the author's edit: the "Summe Monat / Summe Jahr /
# Ueberstunden Monat / Ueberstunden Jahr / Soll Monat / Soll Jahr" mini table
# (previously sitting in F32:J37) and the vacation-tracking table
# (previously A35:C39) are moved further down the sheet (to rows 40-46,
# leaving two blank spacer rows 38-39) so a new "shift date" box (the framed
# G33:I34 area) can be inserted above them. J31 also gets a new bottom/side
# border to close off the J column's data box now that the summary row
# directly below it is gone.
# ============================================================================

# ---- 1. Move the F-J stats block (old rows 32-34) down to rows 40-41, ----
# ---- gaining new G/H columns along the way. ----
$ws.Range("F32:J32").Copy() | Out-Null
$ws.Range("F40").PasteSpecial(-4122) | Out-Null

$ws.Range("F33:J33").Copy() | Out-Null
$ws.Range("F41").PasteSpecial(-4122) | Out-Null

# ---- 2. Move the vacation table (old rows 35-39, cols A-C) plus the ----
# ---- continuing F-J stat rows down to rows 42-46. ----
$ws.Range("A35:C35").Copy() | Out-Null
$ws.Range("A42").PasteSpecial(-4122) | Out-Null
$ws.Range("F33:J33").Copy() | Out-Null
$ws.Range("F42").PasteSpecial(-4122) | Out-Null

$ws.Range("A36:C36").Copy() | Out-Null
$ws.Range("A43").PasteSpecial(-4122) | Out-Null
$ws.Range("F35:J35").Copy() | Out-Null
$ws.Range("F43").PasteSpecial(-4122) | Out-Null

$ws.Range("A37:C37").Copy() | Out-Null
$ws.Range("A44").PasteSpecial(-4122) | Out-Null
$ws.Range("F36:J36").Copy() | Out-Null
$ws.Range("F44").PasteSpecial(-4122) | Out-Null

$ws.Range("A38:C38").Copy() | Out-Null
$ws.Range("A45").PasteSpecial(-4122) | Out-Null
$ws.Range("F37:J37").Copy() | Out-Null
$ws.Range("F45").PasteSpecial(-4122) | Out-Null

$ws.Range("A39:C39").Copy() | Out-Null
$ws.Range("A46").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---- 3. Write the actual values/formulas into their new homes. ----
$ws.Range("I40").Value = "Summe Monat"
$ws.Range("J40").Formula = "=SUM(J4:J31)"

$ws.Range("I41").Value = "Summe Jahr"

$ws.Range("A42").Value = "Urlaubanspruch"
$ws.Range("C42").Value = 30
$ws.Range("I42").Value = "Ueberstunden Monat"

$ws.Range("A43").Value = "Resturlaub Vorjahr"
$ws.Range("I43").Value = "Ueberstunden Jahr"

$ws.Range("A44").Value = "genommen"
$ws.Range("I44").Value = "Soll Monat"

$ws.Range("A45").Value = "Resturlaub "
$ws.Range("I45").Value = "Soll Jahr"

$ws.Range("A46").Value = "Krank "

# ---- 4. New framed "shift date" box: G33 / I34 get a blank, bordered cell ----
$ws.Range("G33").VerticalAlignment = -4108
$ws.Range("I34").ClearContents()
$ws.Range("I34").VerticalAlignment = -4108

# ---- 5. New top border on G40 (closes the bottom of the shift-date box). ----
$ws.Range("G40").VerticalAlignment = -4108
$ws.Range("G40").Borders.Item(8).LineStyle = 1
$ws.Range("G40").Borders.Item(8).Weight = 2
$ws.Range("G40").Borders.Item(8).ColorIndex = 64

# Draw the rest of the shift-date box frame around G33:I34 (thin grey on
# left/right/top, solid border on the inner/outer edges matching the
# original workbook's "light grey box" convention used elsewhere).
$box = $ws.Range("G33:I34")
$box.Borders.Item(7).LineStyle = 1
$box.Borders.Item(7).Weight = 2
$box.Borders.Item(7).ThemeColor = 0
$box.Borders.Item(7).TintAndShade = -0.249977111117893
$box.Borders.Item(10).LineStyle = 1
$box.Borders.Item(10).Weight = 2
$box.Borders.Item(10).ThemeColor = 0
$box.Borders.Item(10).TintAndShade = -0.249977111117893
$box.Borders.Item(8).LineStyle = 1
$box.Borders.Item(8).Weight = 2
$box.Borders.Item(8).ThemeColor = 0
$box.Borders.Item(8).TintAndShade = -0.249977111117893
$ws.Range("G33:I33").Borders.Item(9).LineStyle = 1
$ws.Range("G33:I33").Borders.Item(9).Weight = 2
$ws.Range("G33:I33").Borders.Item(9).ColorIndex = 64

# ---- 6. Update J31's border: close off the bottom of the J-column box ----
# ---- now that the "Summe Monat" row is no longer directly beneath it. ----
$j31 = $ws.Range("J31")
$j31.Borders.Item(7).LineStyle = 1
$j31.Borders.Item(7).Weight = 2
$j31.Borders.Item(7).ThemeColor = 0
$j31.Borders.Item(7).TintAndShade = -0.249977111117893
$j31.Borders.Item(10).LineStyle = 1
$j31.Borders.Item(10).Weight = 2
$j31.Borders.Item(10).ThemeColor = 0
$j31.Borders.Item(10).TintAndShade = -0.249977111117893
$j31.Borders.Item(8).LineStyle = 1
$j31.Borders.Item(8).Weight = 2
$j31.Borders.Item(8).ThemeColor = 0
$j31.Borders.Item(8).TintAndShade = -0.249977111117893
$j31.Borders.Item(9).LineStyle = 1
$j31.Borders.Item(9).Weight = 2
$j31.Borders.Item(9).ColorIndex = 64

# ---- 7. Clear out the cells left behind by the move. ----
$ws.Range("F32").Clear()
$ws.Range("H32").Clear()
$ws.Range("I32").Clear()
$ws.Range("J32").Clear()

$ws.Range("F33").Clear()
$ws.Range("I33").Clear()
$ws.Range("J33").Clear()

$ws.Range("F34").Clear()

$ws.Range("A35").Clear()
$ws.Range("B35").Clear()
$ws.Range("C35").Clear()
$ws.Range("F35").Clear()
$ws.Range("I35").Clear()
$ws.Range("J35").Clear()

$ws.Range("A36").Clear()
$ws.Range("B36").Clear()
$ws.Range("C36").Clear()
$ws.Range("F36").Clear()
$ws.Range("I36").Clear()
$ws.Range("J36").Clear()

$ws.Range("A37").Clear()
$ws.Range("B37").Clear()
$ws.Range("C37").Clear()
$ws.Range("F37").Clear()
$ws.Range("G37").Clear()
$ws.Range("H37").Clear()
$ws.Range("I37").Clear()
$ws.Range("J37").Clear()

$ws.Range("A38").Clear()
$ws.Range("B38").Clear()
$ws.Range("C38").Clear()

$ws.Range("A39").Clear()
$ws.Range("B39").Clear()
$ws.Range("C39").Clear()

# ---- 8. Sheet view bookkeeping: scroll position & selection moved with ----
# ---- the content. ----
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 16 } catch {}
$ws.Range("L40").Select()
